$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: perform all row insertions first (top to bottom), so that the
#     final row numbers below are stable for the subsequent value writes. ---

# Insert new row 8 (becomes "anguilas" / "NA" / "anguila")
$ws.Rows.Item(8).Insert() | Out-Null

# Insert new row before "los demas abalones y orejas de mar"
# (that row was originally #39, now shifted to #40 by the insert above)
$ws.Rows.Item(40).Insert() | Out-Null

# Insert new row before "los demas cangrejos"
# (that row was originally #41, now shifted to #43 by the two inserts above)
$ws.Rows.Item(43).Insert() | Out-Null

# --- Step 2: write cell values in the exact order the new unique strings
#     are first introduced, so the shared-strings table ends up in the
#     same order as the target workbook. ---

# 1) caracol_no_mar (row 16, column C)
$ws.Range("C16").Value = "caracol_no_mar"

# 2) anguila (row 8, column C)
$ws.Range("C8").Value = "anguila"

# 3) anguilas (row 8, columns A and B)
$ws.Range("A8").Value = "anguilas"
$ws.Range("B8").Value = "NA"

# 4) los demas camarones, langostinos y natantia (row 43)
$ws.Range("A43").Value = "los demas camarones, langostinos y natantia"
$ws.Range("B43").Value = "NA"
$ws.Range("C43").Value = "decapodos_otros"

# 5) loco (row 40, columns A and C) then 6) Concholepas concholepas (row 40, column B)
$ws.Range("A40").Value = "loco"
$ws.Range("B40").Value = "Concholepas concholepas"
$ws.Range("C40").Value = "loco"
$ws.Range("A40").WrapText = $true

# --- Step 3: restore the active selection / scroll position to match the
#     edited workbook (best effort - window-position metadata is cosmetic) ---
$ws.Activate()
$ws.Range("C41").Select() | Out-Null
$win = $wb.Windows.Item(1)
$win.ScrollRow = 30
$win.ScrollColumn = 1
